{"js": "const replacements = [\n  [\"2025-02-27 Thursday\", \"2025-02-28 Friday\"],\n  [\"77\u00d797=\", \"97\u00d717=\"],\n  [\"24\u00d757=\", \"92\u00d766=\"],\n  [\"76\u00d743=\", \"98\u00d720=\"],\n  [\"82\u00d797=\", \"65\u00d753=\"],\n  [\"56\u00d770=\", \"70\u00d766=\"],\n  [\"67\u00d756=\", \"15\u00d764=\"],\n  [\"38\u00d795=\", \"78\u00d786=\"],\n  [\"28\u00d797=\", \"31\u00d760=\"],\n  [\"33\u00d740=\", \"83\u00d741=\"],\n  [\"98\u00d780=\", \"99\u00d741=\"],\n  [\"84\u00d761=\", \"38\u00d762=\"],\n  [\"75\u00d799=\", \"21\u00d792=\"],\n  [\"81\u00d744=\", \"98\u00d732=\"],\n  [\"27\u00d787=\", \"91\u00d731=\"],\n  [\"83\u00d744=\", \"16\u00d753=\"],\n  [\"98\u00d796=\", \"13\u00d743=\"],\n  [\"58\u00d745=\", \"52\u00d725=\"],\n  [\"77\u00d735=\", \"69\u00d784=\"],\n  [\"62\u00d799=\", \"13\u00d759=\"],\n  [\"88\u00d716=\", \"48\u00d796=\"],\n  [\"50\u00d718=\", \"41\u00d769=\"],\n  [\"73\u00d738=\", \"14\u00d758=\"],\n  [\"31\u00d722=\", \"40\u00d755=\"],\n  [\"49\u00d733=\", \"27\u00d771=\"],\n  [\"26\u00d713=\", \"19\u00d756=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    @(\"2025-02-27 Thursday\", \"2025-02-28 Friday\"),\n    @(\"77\u00d797=\", \"97\u00d717=\"),\n    @(\"24\u00d757=\", \"92\u00d766=\"),\n    @(\"76\u00d743=\", \"98\u00d720=\"),\n    @(\"82\u00d797=\", \"65\u00d753=\"),\n    @(\"56\u00d770=\", \"70\u00d766=\"),\n    @(\"67\u00d756=\", \"15\u00d764=\"),\n    @(\"38\u00d795=\", \"78\u00d786=\"),\n    @(\"28\u00d797=\", \"31\u00d760=\"),\n    @(\"33\u00d740=\", \"83\u00d741=\"),\n    @(\"98\u00d780=\", \"99\u00d741=\"),\n    @(\"84\u00d761=\", \"38\u00d762=\"),\n    @(\"75\u00d799=\", \"21\u00d792=\"),\n    @(\"81\u00d744=\", \"98\u00d732=\"),\n    @(\"27\u00d787=\", \"91\u00d731=\"),\n    @(\"83\u00d744=\", \"16\u00d753=\"),\n    @(\"98\u00d796=\", \"13\u00d743=\"),\n    @(\"58\u00d745=\", \"52\u00d725=\"),\n    @(\"77\u00d735=\", \"69\u00d784=\"),\n    @(\"62\u00d799=\", \"13\u00d759=\"),\n    @(\"88\u00d716=\", \"48\u00d796=\"),\n    @(\"50\u00d718=\", \"41\u00d769=\"),\n    @(\"73\u00d738=\", \"14\u00d758=\"),\n    @(\"31\u00d722=\", \"40\u00d755=\"),\n    @(\"49\u00d733=\", \"27\u00d771=\"),\n    @(\"26\u00d713=\", \"19\u00d756=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}"}
